# Clear cross elasticities from KDI data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose "Elasticity" (column F) value should be cleared to 0.
# These are exactly the "demand_cross_price" rows (Commodity <> Cross
# commodity) - the "demand_own_price" diagonal rows (F2, F8, F14, F20, F26)
# are left untouched.
$crossRows = 3,4,5,6,7,9,10,11,12,13,15,16,17,18,19,21,22,23,24,25

foreach ($r in $crossRows) {
    $ws.Range("F$r").Value = 0
}

# Add new note text in L14 (introduces a new shared string "elasti" and
# extends the used range of the sheet out to column L).
$ws.Range("L14").Value = "elasti"

# Update the view: scroll back to the top of the sheet and select G28
# instead of L14.
$ws.Activate()
$excel.Goto($ws.Range("D1"), $true)
$ws.Range("G28").Select()
